$d = $word.ActiveDocument

# 1. Update the city line to include the postal code.
$d.Content.Find.Execute("Cluj Napoca", $true, $false, $false, $false, $false, $true, 1, $false, "Cluj Napoca City, 400158", 2) | Out-Null

# 2. Remove the now-redundant "City, State 400158" paragraph entirely (text + paragraph mark).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "City, State 400158") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}
